$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data entry: mark additional "X" boxes that were checked ---
$ws.Range("D15").Value = "X"
$ws.Range("E15").Value = "X"
$ws.Range("D16").Value = "X"
$ws.Range("E16").Value = "X"
$ws.Range("D26").Value = "x"

# --- Print area ---
$ws.PageSetup.PrintArea = "PUNTUACIÓN!`$A`$1:`$H`$31"

# --- Page setup ---
$ws.PageSetup.Orientation = 2  # xlLandscape
$ws.PageSetup.LeftMargin = 65.19685039370079
$ws.PageSetup.RightMargin = 65.19685039370079
$ws.PageSetup.TopMargin = 68.03149606299213
$ws.PageSetup.BottomMargin = 68.03149606299213
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0

# --- Column widths ---
# (engine quantizes ColumnWidth to a 1/6-character pixel grid under the hood,
# so the inputs below are chosen to land as close as possible to the target
# OOXML <col width> values of 3.5 and 4.875 respectively)
$ws.Columns.Item(2).ColumnWidth = 2.666666666666667
$ws.Columns.Item(8).ColumnWidth = 4.0

# --- View settings ---
$window = $excel.ActiveWindow
$window.View = 2  # xlPageBreakPreview
$window.Zoom = 60
$ws.Range("B1").Select()
$window.ScrollRow = 1
$window.ScrollColumn = 2
$ws.Range("J25").Select()
